# Conserto do erro com o rotulo da coluna 2050 nas tabelas e retirada das
# linhas com total das tabelas.
$wb = $excel.ActiveWorkbook

# Helper: replace the (wrongly numeric) E1 header with the correct text
# label, without letting Excel re-interpret a purely-numeric label like
# "2050" back into a number. We briefly mark the cell as Text so the
# assignment sticks as a string, then restore the original "General"
# look-and-feel by re-applying the neighbouring (already textual) D1
# cell's formatting on top, so the style index Excel ends up using for
# E1 matches the one already used by the other header cells.
function Set-YearHeader {
    param($sheetName, $text)

    $ws = $wb.Worksheets.Item($sheetName)
    $ws.Range("E1").NumberFormat = "@"
    $ws.Range("E1").Value = $text
    $ws.Range("D1").Copy()
    $ws.Range("E1").PasteSpecial(-4122)  # xlPasteFormats
}

$excel.CutCopyMode = $false

Set-YearHeader "Potencia Acumulada - SIN (MW)" "2050"
Set-YearHeader "Geracao Periodo Medio (MWMed)" "2050"
Set-YearHeader "Atendimento a Ponta(MW)" "2050"
Set-YearHeader "Potencia Incremental - SIN(MW)" "2041-2050"
Set-YearHeader "Emissoes Totais (MtCO2eq)" "2050"

$excel.CutCopyMode = $false

# Drop the stray "Total" row on every table that has one.
$wb.Worksheets.Item("Potencia Acumulada - SIN (MW)").Rows.Item(13).Delete()
$wb.Worksheets.Item("Geracao Periodo Medio (MWMed)").Rows.Item(13).Delete()
$wb.Worksheets.Item("Atendimento a Ponta(MW)").Rows.Item(13).Delete()
$wb.Worksheets.Item("Potencia Incremental - SIN(MW)").Rows.Item(13).Delete()
$wb.Worksheets.Item("Custo Total (bilhões de R$)").Rows.Item(4).Delete()
